$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (I1, J1) - new columns "I0" and "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from H1 onto the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Data rows 2-11 for columns I and J
$values = @{
    2  = @(1, 2)
    3  = @(1, 6)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 4)
    10 = @(1, 3)
    11 = @(3, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
